$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1000
$ws.Range("D3").Value = 1000
$ws.Range("B5").Value = 5
